$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.370.10'
$ws.Range("E2").Value = '  -0.49%  '

$ws.Range("D3").Value = '3.435.88'
$ws.Range("E3").Value = '  -2.94%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.30'
$ws.Range("E5").Value = '  -2.66%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.56'
$ws.Range("E6").Value = '  -4.72%  '

$ws.Range("D7").Value = '3.434.59'
$ws.Range("E7").Value = '  -2.96%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("E9").Value = '  -4.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.10'
$ws.Range("E11").Value = '  -10.06%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.376'
$ws.Range("E12").Value = '  -7.13%  '

$ws.Range("D13").Value = '4.018.07'
$ws.Range("E13").Value = '  -2.65%  '

$ws.Range("E14").Value = '  -9.29%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '3.506.31'
$ws.Range("E15").Value = '  +0.04%  '

$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.52'
$ws.Range("E16").Value = '  -7.08%  '

$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.115'
$ws.Range("E17").Value = '  -1.93%  '

$ws.Range("D18").Value = '65.207.25'
$ws.Range("E18").Value = '  -0.60%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.69'
$ws.Range("E19").Value = '  -12.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.82'
$ws.Range("E20").Value = '  -6.29%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.53'
$ws.Range("E21").Value = '  -5.45%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '385.55'
$ws.Range("E22").Value = '  -7.65%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.553'
$ws.Range("E23").Value = '  -7.72%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.61'
$ws.Range("E24").Value = '  -6.76%  '

$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.06%  '

$ws.Range("D26").Value = '3.570.57'
$ws.Range("E26").Value = '  -2.67%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000107'
$ws.Range("E27").Value = '  -7.67%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.997'
$ws.Range("E28").Value = '  -0.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.17'
$ws.Range("E29").Value = '  -8.05%  '

$ws.Range("E30").Value = '  -8.46%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.21'
$ws.Range("E31").Value = '  -10.07%  '

$ws.Range("D32").Value = '3.443.33'
$ws.Range("E32").Value = '  -2.95%  '

$ws.Range("E33").Value = '  +0.02%  '

$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.143'
$ws.Range("E34").Value = '  -7.04%  '

$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.16'
$ws.Range("E35").Value = '  -4.83%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '168.72'
$ws.Range("E36").Value = '  -3.52%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.78'
$ws.Range("E37").Value = '  -9.87%  '

$ws.Range("E38").Value = '  -10.49%  '

$ws.Range("E39").Value = '  -7.32%  '

$ws.Range("E40").Value = '  -11.44%  '

$ws.Range("E41").Value = '  -7.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.813'
$ws.Range("E42").Value = '  -5.20%  '

$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.49'
$ws.Range("E43").Value = '  -5.64%  '

$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  +0.06%  '

$ws.Range("E45").Value = '  -14.20%  '

$ws.Range("E46").Value = '  -8.81%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.13'
$ws.Range("E47").Value = '  +2.24%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.78'
$ws.Range("E48").Value = '  -3.12%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.48'
$ws.Range("E49").Value = '  -7.62%  '

$ws.Range("E50").Value = '  -13.08%  '

$ws.Range("D51").Value = '2.166.52'
$ws.Range("E51").Value = '  -6.99%  '
